$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists Modbus registers (Numero / Numero(HEX) / Descripcion).
# A new register "Start a reading if set >0 (cleared by soft)" (register 0x02)
# is being inserted right before the existing "Silo Full" row, so every
# register below it moves down by one row and the numbering (and its
# DEC2HEX hex column) stays contiguous.

# 1) Insert a new row at row 4 (old row 4 "Silo Full" and everything after
#    it shifts down one row).
$ws.Rows(4).Insert()

# 2) The freshly inserted row starts out blank/unformatted - copy the
#    look of the row above (row 3) into it so borders/number formatting
#    match the rest of the table.
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Fill in the new register's data.
$ws.Range("A4").Value = 2
$ws.Range("B4").Formula = "=DEC2HEX(A4,2)"
$ws.Range("C4").Value = "Start a reading if set >0 (cleared by soft)"

# 4) Renumber column A for every row pushed down by the insert (and the
#    brand-new last row) so "Numero" stays a plain 0..76 sequence; column B
#    recalculates automatically since it is a DEC2HEX(A..) formula.
for ($i = 5; $i -le 78; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 2
}

# 5) Leave the selection where the edit was made, matching the saved file.
$null = $ws.Range("C4").Select()
